# Add a new "AmylasePost" data row (row 43) to the protein_info worksheet,
# mirroring the formatting of the other data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("protein_info")

$seqK = "MRFPSIFTAVLFAASSALAAPVNTTTEDETAQIPAEAVIGYSDLEGDFDVAVLPFSNSTNNGLLFINTTIASIAAKEEGVSLDKREEGEPKATPADWRSQSIYFLLTDRFARTDGSTTATCNTADRKYCGGTWQGIIDKLDYIQGMGFTAIWITPVTAQLPQTTAYGDAYHGYWQQDIYSLNENYGTADDLKALSSALHERGMYLMVDVVANHMGYDGAGSSVDYSVFKPFSSQDYFHPFCLIQNYEDQTQVEDCWLGDNTVSLPDLDTTKDVVKNEWYDWVGSLVSNYSIDGLRIDTVKHVQKDFWPGYNKAAGVYCIGEVLDGDPAYTCPYQNVMDGVLNYPIYYPLLNAFKSTSGSMDDLYNMINTVKSDCPDSTLLGTFVENHDNPRFASYTNDIALAKNVAAFIILNDGIPIIYAGQEQHYAGGNDPANREATWLSGYPTDSELYKLIASANAIRNYAISKDTGFVTYKNWPIYKDDTTIAMRKGTDGSQIVTILSNKGASGDSYTLSLSGAGYTAGQQLTEVIGCTTVTVGSDGNVPVPMAGGLPRVLYPTEKLAGSKICSSS"
$seqM = "MRFPSIFTAVLFAASSALAAPVNTTTEDETAQIPAEAVIGYSDLEGDFDVAVLPFSNSTNNGLLFINTTIASIAAKEEGVSLDKREEGEPK"

$row = 43

$ws.Cells.Item($row, 1).Value  = "AmylasePost"
$ws.Cells.Item($row, 2).Value  = "AmylasePost"
$ws.Cells.Item($row, 3).Value  = 1
$ws.Cells.Item($row, 4).Value  = 0
$ws.Cells.Item($row, 5).Value  = 4
$ws.Cells.Item($row, 6).Value  = 1
$ws.Cells.Item($row, 7).Value  = 0
$ws.Cells.Item($row, 8).Value  = 0
$ws.Cells.Item($row, 9).Value  = 0
$ws.Cells.Item($row, 10).Value = "e"
$ws.Cells.Item($row, 11).Value = $seqK
$ws.Cells.Item($row, 12).Value = 569
$ws.Cells.Item($row, 13).Value = $seqM
$ws.Cells.Item($row, 14).Value = 0
$ws.Cells.Item($row, 17).Value = "mingtao pnas"
$ws.Cells.Item($row, 18).Value = "P0C1B3(with another leader sequence)"
$ws.Cells.Item($row, 19).Value = "3 N-glycosylation sites are in leader pro sequence(not account), linker is included in the sp squence(EEGEPK)"

# Match the formatting used by the other data rows (style index 2), leaving
# the unused O/P columns untouched just like the source row.
$ws.Range("A2:N2").Copy() | Out-Null
$ws.Range("A" + $row + ":N" + $row).PasteSpecial(-4122) | Out-Null
$ws.Range("Q2:S2").Copy() | Out-Null
$ws.Range("Q" + $row + ":S" + $row).PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B44").Select() | Out-Null
